$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 'BRICEÑO LUNA , JESSICA ARACELI - 06:30AM - 05:30PM'
$ws.Range("C1").Value = 'MEDINA MARCELO, NAOMI ARIADNA - 05:30PM - 09:15PM'
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""

$ws.Range("B2").Value = 'YANQUI BRAVO, MIRIAN LUZ - 08:30AM - 12:15PM'
$ws.Range("C2").Value = 'SOTELO GONZALES , CAMILA SOFÍA - 01:00PM - 10:00PM'
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""

$ws.Range("B3").Value = 'QUISPE MONDRAGÓN , JUAN ALFONSO - 09:15AM - 01:00PM'
$ws.Range("C3").Value = 'BARRIENTOS JERI, MILAGROS NICOL - 01:45PM - 10:45PM'
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""

$ws.Range("B4").Value = 'MONTEZUMA DEJO, EVELYN BRUNELLA - 07:00AM - 10:45AM'
$ws.Range("C4").Value = 'YOVERA ROBLES, VICTOR EDUARDO - 10:45AM - 02:30PM'
$ws.Range("D4").Value = 'INGA DELGADO, CARLOS DANIEL - 03:00PM - 06:45PM'
$ws.Range("E4").Value = 'BRENIS LÁRTIGA , SEBASTIÁN - 07:00PM - 10:45PM'

$ws.Range("B5").Value = 'MARTINEZ PAZ, ROCIO ESPERANZA - 09:00AM - 06:00PM'
$ws.Range("C5").Value = 'SALAS VILLANUEVA, JAMILA DASHA - 06:00PM - 09:45PM'
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""

$ws.Range("B6").Value = 'POBLETE SAIRE, FIORELLA ESTHER - 07:30AM - 11:15AM'
$ws.Range("C6").Value = 'CARHUARICRA ESPINOZA, FIORELLA NICOLL - 11:15AM - 03:00PM'
$ws.Range("D6").Value = 'RIVERA CARREÑO , DIANA DESIRÉE - 03:30PM - 07:15PM'
$ws.Range("E6").Value = ""

$ws.Range("B7").Value = 'ALVITE CORNEJO, ANGIE LUCERO - 07:30AM - 11:15AM'
$ws.Range("C7").Value = 'CUSI QUISPE, ANDREA ESTEFANY - 11:15AM - 03:00PM'
$ws.Range("D7").Value = 'HUAYNATES ALTAMIRANO, JIM HANS - 04:00PM - 07:45PM'
$ws.Range("E7").Value = ""

$ws.Range("B8").Value = 'MEZA PEREZ, JUAN CRISTOFER - 09:15AM - 01:00PM'
$ws.Range("C8").Value = 'YACILA GRANDEZ, RODRIGO ANDRE - 02:00PM - 05:45PM'
$ws.Range("D8").Value = 'YANAC DAVILA, GERALD RONNY - 05:45PM - 09:30PM'
$ws.Range("E8").Value = ""

$ws.Range("B9").Value = 'MORENO CANCHANYA, ROSMERY - 10:15AM - 02:00PM'
$ws.Range("C9").Value = 'PARICELA TINEO, JAIME DANIEL - 02:00PM - 05:45PM'
$ws.Range("D9").Value = 'MENDOZA DIEGO, ZAIDA VANESSA - 06:00PM - 09:45PM'
$ws.Range("E9").Value = ""

$ws.Range("B10").Value = 'ILDEFONSO MOTTA, JHOSSEP ANGELO - 10:15AM - 02:00PM'
$ws.Range("C10").Value = 'VARGAS CASTRO, LOANA VICTORIA - 02:00PM - 05:45PM'
$ws.Range("D10").Value = 'FLORES PAREDES, LOURDES - 06:00PM - 08:45PM'
$ws.Range("E10").Value = ""

$ws.Range("B11").Value = 'ZAVALA SOSA, NICOLE - 11:00AM - 08:00PM'
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("B12").Value = 'AYALA MORA, CECILIA ROSARIO - 10:30AM - 02:15PM'
$ws.Range("C12").Value = 'VILCHEZ CUBA, JACK ANTHONY - 04:15PM - 08:00PM'
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("B13").Value = 'GOMEZ ALBINO, IDALIA GIMENA - 11:15AM - 03:00PM'
$ws.Range("C13").Value = 'ARIAS MACHACUAY, SADELITH SORAGGI - 04:30PM - 08:15PM'
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("B14").Value = 'HUAMAN HUAMANI, ALEXIS JAVIER - 11:15AM - 03:00PM'
$ws.Range("C14").Value = 'MENDOZA CRUZ, LILIANA LILIANA - 04:45PM - 08:30PM'
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("B15").Value = 'HUAYANAY VELASCO, ATHINA - 11:15AM - 03:00PM'
$ws.Range("C15").Value = 'TORRES RAZURI, JESUS GUSTAVO SANTIAGO - 05:30PM - 09:15PM'
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
